# Weekly refresh of the "Femacal de La Calera - Espárragos" price log.
# A new weekly reading is inserted at row 8 (pushing the existing rows
# 8-56 down to 9-57), and the new row is populated with this week's
# market data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 8, shifting rows 8:56 down to 9:57.
$ws.Rows("8:8").Insert()

# Populate the newly inserted row 8 with the latest reading.
$newRow = @(
    3,                        # A: Mercado ID
    "Femacal de La Calera",   # B: Mercado
    "Coquimbo",                # C: Región
    44901,                      # D: Fecha (2022-12-06)
    5,                          # E: Codreg
    300000000,                  # F: Categoría ID
    "Espárragos",               # G: Categoría
    "Verde",                    # H: Variedad
    "Primera",                  # I: Calidad
    830,                        # J: Volumen
    1400,                       # K: Precio mínimo
    1500,                       # L: Precio máximo
    1451,                       # M: Precio promedio ponderado
    "`$/kilo",                  # N: Unidad de comercialización
    "Provincia de Quillota",    # O: Origen
    1451,                       # P: Precio $/Kg
    1,                          # Q: Kg o Unidades
    "Hortaliza"                 # R: Clasificación
)

for ($i = 0; $i -lt $newRow.Length; $i++) {
    $ws.Cells.Item(8, $i + 1).Value = $newRow[$i]
}
